$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7-12 first (they are no longer part of the dataset)
$ws.Range("A7:H12").Delete()

# Update row 1
$ws.Range("A1").Value = "0-C-0"
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 0.6
$ws.Range("D1").Value = 0.01112081887113899

# Update row 2
$ws.Range("A2").Value = "1-C-0"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.6
$ws.Range("D2").Value = 0.01112081887113899

# Update row 3
$ws.Range("A3").Value = "2-C-0"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.6
$ws.Range("D3").Value = 0.01112081887113899

# Update row 4
$ws.Range("A4").Value = "3-C-0"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0.6
$ws.Range("D4").Value = 0.01112081887113899

# Update row 5
$ws.Range("A5").Value = "4-C-0"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0.6
$ws.Range("D5").Value = 0.01112081887113899

# Update row 6
$ws.Range("A6").Value = "5-C-0"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 0.6
$ws.Range("D6").Value = 0.01112081887113899
